# Auto-generated script to apply market-data refresh edits to Diabolos_Profits workbook
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has columns H-N holding
# refreshed marketboard pricing data (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 16689158
$ws.Range("I86").Value = 9983.666999999999
$ws.Range("K86").Value = 9983.666999999999
$ws.Range("M86").Value = -8860.666999999999
$ws.Range("H88").Value = 1038
$ws.Range("I88").Value = 943.375
$ws.Range("K88").Value = 943.375
$ws.Range("M88").Value = -537.375
$ws.Range("H89").Value = 16689158
$ws.Range("I89").Value = 9983.666999999999
$ws.Range("K89").Value = 49918.335
$ws.Range("M89").Value = -44302.335
$ws.Range("H91").Value = 1038
$ws.Range("I91").Value = 943.375
$ws.Range("K91").Value = 943.375
$ws.Range("M91").Value = 460.625
$ws.Range("H97").Value = 992
$ws.Range("J97").Value = 992
$ws.Range("L97").Value = 2976
$ws.Range("N97").Value = -3968
$ws.Range("H113").Value = 45458340
$ws.Range("I113").Value = 76926300
$ws.Range("J113").Value = 4611.222
$ws.Range("K113").Value = 76926300
$ws.Range("L113").Value = 4611.222
$ws.Range("M113").Value = -76923046
$ws.Range("N113").Value = -11119.222
$ws.Range("H118").Value = 190
$ws.Range("I118").Value = 190
$ws.Range("K118").Value = 570
$ws.Range("M118").Value = 1087
$ws.Range("H132").Value = 4785.8887
$ws.Range("I132").Value = 4488.909
$ws.Range("J132").Value = 6092.6
$ws.Range("K132").Value = 13466.727
$ws.Range("L132").Value = 18277.8
$ws.Range("M132").Value = -10936.727
$ws.Range("N132").Value = -23337.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 33335152
$ws.Range("I61").Value = 37038836
$ws.Range("K61").Value = 37038836
$ws.Range("M61").Value = -37038624
$ws.Range("H74").Value = 2536.8076
$ws.Range("I74").Value = 1717.762
$ws.Range("J74").Value = 5976.8
$ws.Range("K74").Value = 1717.762
$ws.Range("L74").Value = 5976.8
$ws.Range("M74").Value = -843.7619999999999
$ws.Range("N74").Value = -7724.8
$ws.Range("H77").Value = 2536.8076
$ws.Range("I77").Value = 1717.762
$ws.Range("J77").Value = 5976.8
$ws.Range("K77").Value = 8588.809999999999
$ws.Range("L77").Value = 29884
$ws.Range("M77").Value = -4220.809999999999
$ws.Range("N77").Value = -38620
$ws.Range("H97").Value = 774.64514
$ws.Range("I97").Value = 579.05
$ws.Range("J97").Value = 1130.2727
$ws.Range("K97").Value = 579.05
$ws.Range("L97").Value = 1130.2727
$ws.Range("M97").Value = -83.04999999999995
$ws.Range("N97").Value = -2122.2727
$ws.Range("H132").Value = 43480876
$ws.Range("I132").Value = 55557790
$ws.Range("J132").Value = 3969.4
$ws.Range("K132").Value = 166673370
$ws.Range("L132").Value = 11908.2
$ws.Range("M132").Value = -166670840
$ws.Range("N132").Value = -16968.2
$ws.Range("H136").Value = 33335152
$ws.Range("I136").Value = 37038836
$ws.Range("K136").Value = 111116508
$ws.Range("M136").Value = -111113958

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 145
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 1020.1818
$ws.Range("I99").Value = 876.5
$ws.Range("K99").Value = 876.5
$ws.Range("M99").Value = 621.5
$ws.Range("H107").Value = 27650.158
$ws.Range("I107").Value = 20529.533
$ws.Range("J107").Value = 54352.5
$ws.Range("K107").Value = 20529.533
$ws.Range("L107").Value = 54352.5
$ws.Range("M107").Value = -18609.533
$ws.Range("N107").Value = -58192.5
$ws.Range("H134").Value = 1686.2667
$ws.Range("I134").Value = 1444.975
$ws.Range("K134").Value = 4334.924999999999
$ws.Range("M134").Value = -1799.924999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19106.428
$ws.Range("I41").Value = 3248
$ws.Range("J41").Value = 31000.25
$ws.Range("K41").Value = 3248
$ws.Range("L41").Value = 31000.25
$ws.Range("M41").Value = -2820
$ws.Range("N41").Value = -31856.25
$ws.Range("H50").Value = 35211.145
$ws.Range("J50").Value = 36496.332
$ws.Range("L50").Value = 36496.332
$ws.Range("N50").Value = -37746.332
$ws.Range("H51").Value = 59499.625
$ws.Range("I51").Value = 58999.668
$ws.Range("J51").Value = 60999.5
$ws.Range("K51").Value = 58999.668
$ws.Range("L51").Value = 60999.5
$ws.Range("M51").Value = -58263.668
$ws.Range("N51").Value = -62471.5
$ws.Range("H59").Value = 44999.5
$ws.Range("I59").Value = 40000
$ws.Range("J59").Value = 49999
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 49999
$ws.Range("M59").Value = -38855
$ws.Range("N59").Value = -52289
$ws.Range("H60").Value = 22400
$ws.Range("J60").Value = 24000
$ws.Range("L60").Value = 24000
$ws.Range("N60").Value = -25022
$ws.Range("H61").Value = 59499.625
$ws.Range("I61").Value = 58999.668
$ws.Range("J61").Value = 60999.5
$ws.Range("K61").Value = 58999.668
$ws.Range("L61").Value = 60999.5
$ws.Range("M61").Value = -58651.668
$ws.Range("N61").Value = -61695.5
$ws.Range("H68").Value = 69996.25
$ws.Range("J68").Value = 69996.25
$ws.Range("L68").Value = 69996.25
$ws.Range("N68").Value = -71494.25
$ws.Range("H71").Value = 69996.25
$ws.Range("J71").Value = 69996.25
$ws.Range("L71").Value = 209988.75
$ws.Range("N71").Value = -217476.75
$ws.Range("H132").Value = 2677.04
$ws.Range("I132").Value = 2580.2917
$ws.Range("K132").Value = 7740.875100000001
$ws.Range("M132").Value = -5210.875100000001
$ws.Range("H134").Value = 2924.375
$ws.Range("I134").Value = 2566.3333
$ws.Range("K134").Value = 7698.999899999999
$ws.Range("M134").Value = -5163.999899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14706473
$ws.Range("J2").Value = 17857860
$ws.Range("L2").Value = 107147160
$ws.Range("N2").Value = -107147386
$ws.Range("H112").Value = 1766.3334
$ws.Range("I112").Value = 1766.3334
$ws.Range("K112").Value = 5299.0002
$ws.Range("M112").Value = -4191.0002
$ws.Range("H120").Value = 32716.1
$ws.Range("I120").Value = 25399.2
$ws.Range("K120").Value = 76197.60000000001
$ws.Range("M120").Value = -71359.60000000001
$ws.Range("H139").Value = 62501588
$ws.Range("I139").Value = 83334630
$ws.Range("J139").Value = 2449.5
$ws.Range("K139").Value = 250003890
$ws.Range("L139").Value = 7348.5
$ws.Range("M139").Value = -249998750
$ws.Range("N139").Value = -17628.5
$ws.Range("H140").Value = 1973.9333
$ws.Range("I140").Value = 1291.7273
$ws.Range("J140").Value = 3850
$ws.Range("K140").Value = 3875.1819
$ws.Range("L140").Value = 11550
$ws.Range("M140").Value = 1304.8181
$ws.Range("N140").Value = -21910
$ws.Range("H141").Value = 2909.5
$ws.Range("I141").Value = 2911.4
$ws.Range("J141").Value = 2900
$ws.Range("K141").Value = 8734.200000000001
$ws.Range("L141").Value = 8700
$ws.Range("M141").Value = -3554.200000000001
$ws.Range("N141").Value = -19060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3807.025
$ws.Range("I132").Value = 3051.84
$ws.Range("J132").Value = 5065.6665
$ws.Range("K132").Value = 9155.52
$ws.Range("L132").Value = 15196.9995
$ws.Range("M132").Value = -6625.52
$ws.Range("N132").Value = -20256.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2996.3333
$ws.Range("I40").Value = 2996.3333
$ws.Range("K40").Value = 2996.3333
$ws.Range("M40").Value = -2860.3333
$ws.Range("H46").Value = 2533.7
$ws.Range("I46").Value = 968.5
$ws.Range("J46").Value = 2925
$ws.Range("K46").Value = 968.5
$ws.Range("L46").Value = 2925
$ws.Range("M46").Value = -780.5
$ws.Range("N46").Value = -3301
$ws.Range("H61").Value = 16663.834
$ws.Range("I61").Value = 10512.75
$ws.Range("J61").Value = 28966
$ws.Range("K61").Value = 10512.75
$ws.Range("L61").Value = 28966
$ws.Range("M61").Value = -10310.75
$ws.Range("N61").Value = -29370
$ws.Range("H68").Value = 44001.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 44001.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 44001.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -45499.5
$ws.Range("H71").Value = 44001.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 44001.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 220007.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -227495.5
$ws.Range("H113").Value = 16663.834
$ws.Range("I113").Value = 10512.75
$ws.Range("J113").Value = 28966
$ws.Range("K113").Value = 10512.75
$ws.Range("L113").Value = 28966
$ws.Range("M113").Value = -8342.75
$ws.Range("N113").Value = -33306
$ws.Range("H136").Value = 2252.6943
$ws.Range("I136").Value = 2051.862
$ws.Range("K136").Value = 6155.586
$ws.Range("M136").Value = -3605.586

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3973500.5
$ws.Range("I62").Value = 4766400.5
$ws.Range("K62").Value = 4766400.5
$ws.Range("M62").Value = -4765776.5
$ws.Range("H65").Value = 3973500.5
$ws.Range("I65").Value = 4766400.5
$ws.Range("K65").Value = 23832002.5
$ws.Range("M65").Value = -23828882.5
$ws.Range("H113").Value = 385.44
$ws.Range("I113").Value = 247.85715
$ws.Range("K113").Value = 743.5714499999999
$ws.Range("M113").Value = 1426.42855
$ws.Range("H136").Value = 2394
$ws.Range("I136").Value = 1129.7727
$ws.Range("K136").Value = 3389.3181
$ws.Range("M136").Value = -839.3181
